$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 2: only B2 changes
$ws.Range("B2").Value = 95707

# Rows 3 and 5 swap their A, E, F, G, H, Q, R values, and B changes independently.
# Capture original row 3 values (which become the new row 5 values).
$A3 = $ws.Range("A3").Value2
$E3 = $ws.Range("E3").Value2
$F3 = $ws.Range("F3").Value2
$G3 = $ws.Range("G3").Value2
$H3 = $ws.Range("H3").Value2
$Q3 = $ws.Range("Q3").Value2
$R3 = $ws.Range("R3").Value2

# Capture original row 5 values (which become the new row 3 values).
$A5 = $ws.Range("A5").Value2
$E5 = $ws.Range("E5").Value2
$F5 = $ws.Range("F5").Value2
$G5 = $ws.Range("G5").Value2
$H5 = $ws.Range("H5").Value2
$Q5 = $ws.Range("Q5").Value2
$R5 = $ws.Range("R5").Value2

# Row 3 gets the former row 5 data
$ws.Range("A3").Value = $A5
$ws.Range("B3").Value = 77053
$ws.Range("E3").Value = $E5
$ws.Range("F3").Value = $F5
$ws.Range("G3").Value = $G5
$ws.Range("H3").Value = $H5
$ws.Range("Q3").Value = $Q5
$ws.Range("R3").Value = $R5

# Row 4: only B4 changes
$ws.Range("B4").Value = 77650

# Row 5 gets the former row 3 data
$ws.Range("A5").Value = $A3
$ws.Range("B5").Value = 77402
$ws.Range("E5").Value = $E3
$ws.Range("F5").Value = $F3
$ws.Range("G5").Value = $G3
$ws.Range("H5").Value = $H3
$ws.Range("Q5").Value = $Q3
$ws.Range("R5").Value = $R3

# Row 6: only B6 changes
$ws.Range("B6").Value = 90826
